$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (Text number format) for price cells whose new values look like plain
# decimal numbers, so Excel keeps storing them as text strings (matching the source
# data which uses inline/shared strings for the Price column), instead of silently
# converting them to numeric cell values.
$textCells = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D11",
    "D14",
    "D15",
    "D17",
    "D20",
    "D21",
    "D23",
    "D26",
    "D27",
    "D29",
    "D30",
    "D33",
    "D34",
    "D36",
    "D37",
    "D38",
    "D40",
    "D41",
    "D42",
    "D45",
    "D46",
    "D47",
    "D49",
    "D50",
    "D51",
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated cell values scraped for this run.
$ws.Range("D2").Value = "40.909.12"
$ws.Range("E2").Value = "  -1.86%  "
$ws.Range("D3").Value = "2.407.51"
$ws.Range("E3").Value = "  -2.74%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "313.76"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").Value = "87.81"
$ws.Range("E6").Value = "  -5.75%  "
$ws.Range("D7").Value = "0.533"
$ws.Range("E7").Value = "  -3.47%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("D9").Value = "0.492"
$ws.Range("E9").Value = "  -4.86%  "
$ws.Range("D10").Value = "0.0827"
$ws.Range("E10").Value = "  -3.76%  "
$ws.Range("D11").Value = "31.11"
$ws.Range("E11").Value = "  -6.26%  "
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("D13").Value = "2.780.68"
$ws.Range("E13").Value = "  -2.65%  "
$ws.Range("D14").Value = "6.69"
$ws.Range("E14").Value = "  -3.05%  "
$ws.Range("D15").Value = "15.33"
$ws.Range("E15").Value = "  -2.89%  "
$ws.Range("D16").Value = "2.399.21"
$ws.Range("E16").Value = "  -2.87%  "
$ws.Range("D17").Value = "0.764"
$ws.Range("E17").Value = "  -3.18%  "
$ws.Range("D18").Value = "40.850.61"
$ws.Range("E18").Value = "  -1.90%  "
$ws.Range("D19").Value = "0.0₃0916"
$ws.Range("E19").Value = "  -3.86%  "
$ws.Range("D20").Value = "6.18"
$ws.Range("E20").Value = "  -4.66%  "
$ws.Range("D21").Value = "70.29"
$ws.Range("E21").Value = "  -1.39%  "
$ws.Range("E22").Value = "  -4.63%  "
$ws.Range("D23").Value = "238.30"
$ws.Range("E23").Value = "  -0.59%  "
$ws.Range("E24").Value = "  -3.64%  "
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "1.83"
$ws.Range("E26").Value = "  -5.34%  "
$ws.Range("D27").Value = "23.79"
$ws.Range("E27").Value = "  -4.06%  "
$ws.Range("E28").Value = "  -2.74%  "
$ws.Range("D29").Value = "9.44"
$ws.Range("E29").Value = "  -4.21%  "
$ws.Range("D30").Value = "33.93"
$ws.Range("E30").Value = "  -6.16%  "
$ws.Range("E31").Value = "  -1.08%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("D33").Value = "5.24"
$ws.Range("E33").Value = "  -5.39%  "
$ws.Range("D34").Value = "0.0733"
$ws.Range("E34").Value = "  -4.61%  "
$ws.Range("E35").Value = "  -5.13%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").Value = "0.114"
$ws.Range("E36").Value = "  -1.94%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "2.83"
$ws.Range("E37").Value = "  -3.61%  "
$ws.Range("D38").Value = "16.02"
$ws.Range("E38").Value = "  -7.54%  "
$ws.Range("E39").Value = "  -7.48%  "
$ws.Range("D40").Value = "0.0990"
$ws.Range("E40").Value = "  -4.73%  "
$ws.Range("D41").Value = "3.83"
$ws.Range("E41").Value = "  -4.55%  "
$ws.Range("D42").Value = "2.28"
$ws.Range("E42").Value = "  -7.07%  "
$ws.Range("D43").Value = "1.981.82"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("E44").Value = "  -4.80%  "
$ws.Range("D45").Value = "17.94"
$ws.Range("E45").Value = "  -5.77%  "
$ws.Range("D46").Value = "2.82"
$ws.Range("E46").Value = "  -5.73%  "
$ws.Range("D47").Value = "9.27"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("D48").Value = "2.648.36"
$ws.Range("E48").Value = "  -2.37%  "
$ws.Range("D49").Value = "73.42"
$ws.Range("E49").Value = "  -1.26%  "
$ws.Range("D50").Value = "93.30"
$ws.Range("E50").Value = "  -4.21%  "
$ws.Range("D51").Value = "50.90"
$ws.Range("E51").Value = "  -2.91%  "
